# Update grain orientation logic to use color codes.
# Rows 2-18 are updated in place to their new values (name, code, width,
# height, thickness, color, qty, grain). Rows 19-25 are removed entirely,
# shrinking the used range from A1:H25 down to A1:H18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowNums = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18)
$colA = @("柜体侧板（R）", "柜体侧板（L）", "顶底板（双门柜体）", "顶底板（双门柜体）", "层隔板（双门柜体）", "层隔板（双门柜体）", "层隔板（单门柜体）", "后背板（双门柜体）为(1+1)组合", "门板（L/R）", "门板（L/R）", "门板（L/R）", "门板（L/R）", "单抽屉组件", "底支撑(双门柜体）", "收口条", "底支撑(双门柜体）", "底支撑(双门柜体）")
$colB = @("CB(R)-HS00-2434-574-16", "CB(L)-HS00-2434-574-16", "DD-HS00-736-554-16", "DD-HS00-768-554-16", "CG-HS00-864-554-25", "CG-HS00-800-554-25", "CG-HS00-336-554-25", "HB-HS00-2320-432-12", "MB(R)-(门板花色)-2352-365-16", "MB(L)-(门板花色)-2288-429-16", "MB(R)-(门板花色)-2320-349-16", "MB(L)-(门板花色)-2320-381-16", "抽屉拉板-HS00-438-106-12", "DC-HS98-704-82-16", "TSB50-HS03-2434-50-16", "DC-HS98-896-82-16", "DC-HS98-672-82-16")
$colC = @(2434, 2434, 736, 768, 864, 800, 336, 2320, 2352, 2288, 2320, 2320, 438, 704, 2434, 896, 672)
$colD = @(574, 574, 554, 554, 554, 554, 554, 432, 365, 429, 349, 381, 106, 82, 50, 82, 82)
$colE = @(16, 16, 16, 16, 25, 25, 25, 12, 16, 16, 16, 16, 12, 16, 16, 16, 16)
$colF = @("HS00", "HS00", "HS00", "HS00", "HS00", "HS00", "HS00", "HS00", "(门板花色)", "(门板花色)", "(门板花色)", "(门板花色)", "HS00", "HS98", "HS03", "HS98", "HS98")
$colG = @(1, 1, 1, 1, 1, 3, 2, 1, 1, 1, 1, 2, 3, 1, 2, 2, 2)
$colH = @("mixed", "mixed", "mixed", "mixed", "mixed", "mixed", "mixed", "mixed", "fixed", "fixed", "fixed", "fixed", "mixed", "mixed", "fixed", "mixed", "mixed")

for ($i = 0; $i -lt $rowNums.Count; $i++) {
    $r = $rowNums[$i]
    $ws.Cells.Item($r, 1).Value = $colA[$i]
    $ws.Cells.Item($r, 2).Value = $colB[$i]
    $ws.Cells.Item($r, 3).Value = $colC[$i]
    $ws.Cells.Item($r, 4).Value = $colD[$i]
    $ws.Cells.Item($r, 5).Value = $colE[$i]
    $ws.Cells.Item($r, 6).Value = $colF[$i]
    $ws.Cells.Item($r, 7).Value = $colG[$i]
    $ws.Cells.Item($r, 8).Value = $colH[$i]
}

# Remove the now-obsolete trailing rows (19-25); this also shrinks the
# sheet's used range/dimension down to A1:H18 automatically.
$ws.Range("A19:H25").Delete()
